$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# 1) Remove the now-unwanted trailing " " run from the first paragraph
#    (do this before any text-length-changing Find/Replace so offsets stay valid).
$delRange = $d.Range(31, 32)
$delRange.Delete()

# 2) Replace the placeholder text in the remaining run.
$d.Content.Find.Execute("**ID__AFFARS_5310_topic_1__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_PART_5310__ID**", 2)

# 3) Update the paragraph's left indent: 120 twips (6pt) -> 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# 4) Add a paragraph border (top/left/bottom/right) with 5pt distance from text,
#    matching <w:pBdr><w:top w:space="5"/>...</w:pBdr> (no line style/size/color).
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromRight = 5
